$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold/border/center-top alignment) from the existing
# date column (A2) down to the newly added rows (A50:A68) before writing data,
# so the new cells inherit the same style used throughout column A.
$ws.Range("A2").Copy($ws.Range("A50:A68"))

# Rewrite the full data block (rows 2-68) in the new row order, which
# interleaves each year's Oct/Nov/Dec ahead of Jan-Sep, and appends the
# new 2022 and partial-2023 rows.
$ws.Range("A2").Value = "2018-10"
$ws.Range("B2").Value = 105.4
$ws.Range("C2").Value = 105.4
$ws.Range("A3").Value = "2018-11"
$ws.Range("B3").Value = 104.7
$ws.Range("C3").Value = 104.7
$ws.Range("A4").Value = "2018-12"
$ws.Range("B4").Value = 103.9
$ws.Range("C4").Value = 103.9
$ws.Range("A5").Value = "2018-01"
$ws.Range("B5").Value = 104.516
$ws.Range("C5").Value = 104.516
$ws.Range("A6").Value = "2018-02"
$ws.Range("B6").Value = 103.2
$ws.Range("C6").Value = 103.2
$ws.Range("A7").Value = "2018-03"
$ws.Range("B7").Value = 102.2
$ws.Range("C7").Value = 102.2
$ws.Range("A8").Value = "2018-04"
$ws.Range("B8").Value = 101.1
$ws.Range("C8").Value = 101.1
$ws.Range("A9").Value = "2018-05"
$ws.Range("B9").Value = 103.1
$ws.Range("C9").Value = 103.1
$ws.Range("A10").Value = "2018-06"
$ws.Range("B10").Value = 104.3
$ws.Range("C10").Value = 104.3
$ws.Range("A11").Value = "2018-07"
$ws.Range("B11").Value = 105
$ws.Range("C11").Value = 105
$ws.Range("A12").Value = "2018-08"
$ws.Range("B12").Value = 105.2
$ws.Range("C12").Value = 105.2
$ws.Range("A13").Value = "2018-09"
$ws.Range("B13").Value = 105.4
$ws.Range("C13").Value = 105.4
$ws.Range("A14").Value = "2019-10"
$ws.Range("B14").Value = 99.7
$ws.Range("C14").Value = 99.7
$ws.Range("A15").Value = "2019-11"
$ws.Range("B15").Value = 100.4
$ws.Range("C15").Value = 100.4
$ws.Range("A16").Value = "2019-12"
$ws.Range("B16").Value = 99.7
$ws.Range("C16").Value = 99.7
$ws.Range("A17").Value = "2019-01"
$ws.Range("B17").Value = 104.9
$ws.Range("C17").Value = 104.9
$ws.Range("A18").Value = "2019-02"
$ws.Range("B18").Value = 105.3
$ws.Range("C18").Value = 105.3
$ws.Range("A19").Value = "2019-03"
$ws.Range("B19").Value = 106.8
$ws.Range("C19").Value = 106.8
$ws.Range("A20").Value = "2019-04"
$ws.Range("B20").Value = 106.5
$ws.Range("C20").Value = 106.5
$ws.Range("A21").Value = "2019-05"
$ws.Range("B21").Value = 104.9
$ws.Range("C21").Value = 104.9
$ws.Range("A22").Value = "2019-06"
$ws.Range("B22").Value = 104.4
$ws.Range("C22").Value = 104.4
$ws.Range("A23").Value = "2019-07"
$ws.Range("B23").Value = 103
$ws.Range("C23").Value = 103
$ws.Range("A24").Value = "2019-08"
$ws.Range("B24").Value = 101.3
$ws.Range("C24").Value = 101.3
$ws.Range("A25").Value = "2019-09"
$ws.Range("B25").Value = 100.3
$ws.Range("C25").Value = 100.3
$ws.Range("A26").Value = "2020-10"
$ws.Range("B26").Value = 94.59999999999999
$ws.Range("C26").Value = 94.59999999999999
$ws.Range("A27").Value = "2020-11"
$ws.Range("B27").Value = 96.09999999999999
$ws.Range("C27").Value = 96.09999999999999
$ws.Range("A28").Value = "2020-12"
$ws.Range("B28").Value = 97.90000000000001
$ws.Range("C28").Value = 97.90000000000001
$ws.Range("A29").Value = "2020-01"
$ws.Range("B29").Value = 99.40000000000001
$ws.Range("C29").Value = 99.40000000000001
$ws.Range("A30").Value = "2020-02"
$ws.Range("B30").Value = 98.7
$ws.Range("C30").Value = 98.7
$ws.Range("A31").Value = "2020-03"
$ws.Range("B31").Value = 96.40000000000001
$ws.Range("C31").Value = 96.40000000000001
$ws.Range("A32").Value = "2020-04"
$ws.Range("B32").Value = 94.09999999999999
$ws.Range("C32").Value = 94.09999999999999
$ws.Range("A33").Value = "2020-05"
$ws.Range("B33").Value = 94.2
$ws.Range("C33").Value = 94.2
$ws.Range("A34").Value = "2020-06"
$ws.Range("B34").Value = 93.5
$ws.Range("C34").Value = 93.5
$ws.Range("A35").Value = "2020-07"
$ws.Range("B35").Value = 94.3
$ws.Range("C35").Value = 94.3
$ws.Range("A36").Value = "2020-08"
$ws.Range("B36").Value = 94.7
$ws.Range("C36").Value = 94.7
$ws.Range("A37").Value = "2020-09"
$ws.Range("B37").Value = 94.5
$ws.Range("C37").Value = 94.5
$ws.Range("A38").Value = "2021-10"
$ws.Range("B38").Value = 109.4
$ws.Range("C38").Value = 99.59999999999999
$ws.Range("A39").Value = "2021-11"
$ws.Range("B39").Value = 110.9
$ws.Range("C39").Value = 99.8
$ws.Range("A40").Value = "2021-12"
$ws.Range("B40").Value = 112.5
$ws.Range("C40").Value = 100
$ws.Range("A41").Value = "2021-01"
$ws.Range("B41").Value = 100.3
$ws.Range("C41").Value = 91.5
$ws.Range("A42").Value = "2021-02"
$ws.Range("B42").Value = 100.8
$ws.Range("C42").Value = 91.5
$ws.Range("A43").Value = "2021-03"
$ws.Range("B43").Value = 101.8
$ws.Range("C43").Value = 91.8
$ws.Range("A44").Value = "2021-04"
$ws.Range("B44").Value = 102.4
$ws.Range("C44").Value = 92.5
$ws.Range("A45").Value = "2021-05"
$ws.Range("B45").Value = 102.2
$ws.Range("C45").Value = 92.5
$ws.Range("A46").Value = "2021-06"
$ws.Range("B46").Value = 103.5
$ws.Range("C46").Value = 92.5
$ws.Range("A47").Value = "2021-07"
$ws.Range("B47").Value = 103.6
$ws.Range("C47").Value = 92.2
$ws.Range("A48").Value = "2021-08"
$ws.Range("B48").Value = 105.6
$ws.Range("C48").Value = 98.5
$ws.Range("A49").Value = "2021-09"
$ws.Range("B49").Value = 108.5
$ws.Range("C49").Value = 99.59999999999999
$ws.Range("A50").Value = "2022-10"
$ws.Range("B50").Value = 115.3
$ws.Range("C50").Value = 100
$ws.Range("A51").Value = "2022-11"
$ws.Range("B51").Value = 112.5
$ws.Range("C51").Value = 100
$ws.Range("A52").Value = "2022-12"
$ws.Range("B52").Value = 108.4
$ws.Range("C52").Value = 100
$ws.Range("A53").Value = "2022-01"
$ws.Range("B53").Value = 109.9
$ws.Range("C53").Value = 100
$ws.Range("A54").Value = "2022-02"
$ws.Range("B54").Value = 110.4
$ws.Range("C54").Value = 100
$ws.Range("A55").Value = "2022-03"
$ws.Range("B55").Value = 116.3
$ws.Range("C55").Value = 100
$ws.Range("A56").Value = "2022-04"
$ws.Range("B56").Value = 120.7
$ws.Range("C56").Value = 100
$ws.Range("A57").Value = "2022-05"
$ws.Range("B57").Value = 122
$ws.Range("C57").Value = 100
$ws.Range("A58").Value = "2022-06"
$ws.Range("B58").Value = 121.8
$ws.Range("C58").Value = 100
$ws.Range("A59").Value = "2022-07"
$ws.Range("B59").Value = 120.7
$ws.Range("C59").Value = 100
$ws.Range("A60").Value = "2022-08"
$ws.Range("B60").Value = 118.6
$ws.Range("C60").Value = 100
$ws.Range("A61").Value = "2022-09"
$ws.Range("B61").Value = 116.4
$ws.Range("C61").Value = 100
$ws.Range("A62").Value = "2023-01"
$ws.Range("B62").Value = 108
$ws.Range("C62").Value = 100
$ws.Range("A63").Value = "2023-02"
$ws.Range("B63").Value = 109
$ws.Range("C63").Value = 100
$ws.Range("A64").Value = "2023-03"
$ws.Range("B64").Value = 104.5
$ws.Range("C64").Value = 100
$ws.Range("A65").Value = "2023-04"
$ws.Range("B65").Value = 101.4
$ws.Range("C65").Value = 100
$ws.Range("A66").Value = "2023-05"
$ws.Range("B66").Value = 99.59999999999999
$ws.Range("C66").Value = 100
$ws.Range("A67").Value = "2023-06"
$ws.Range("B67").Value = 98.5
$ws.Range("C67").Value = 100
$ws.Range("A68").Value = "2023-07"
$ws.Range("B68").Value = 98.8
$ws.Range("C68").Value = 100
